$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at H:I (pushes existing attackOrigin/splashType/etc. right by 2)
$null = $ws.Columns("H:I").Insert()

# New column headers/values: targetAllies (bool) and targetEnemies (bool)
$ws.Range("H1").Value = "targetAllies"
$ws.Range("I1").Value = "targetEnemies"
$ws.Range("H2").Value = $false
$ws.Range("I2").Value = $true

# Rename the (now shifted) attackOrigin header to splashOrigin; its value cell is untouched
$ws.Range("J1").Value = "splashOrigin"

# Tweak a couple of stat values while working on splash attack functionality
$ws.Range("BL2").Value = 3
$ws.Range("CN2").Value = 0.1

# Restore selection to the new targetEnemies cell
$null = $ws.Range("I2").Select()
